$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.690.39"
$c.Style = $savedStyle
$ws.Range("E2").Value = "  -0.01%  "
$c = $ws.Range("D3")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.634.54"
$c.Style = $savedStyle
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "212.12"
$c.Style = $savedStyle
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.40%  "
$c = $ws.Range("D7")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $savedStyle
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "23.36"
$c.Style = $savedStyle
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  +0.22%  "
$c = $ws.Range("D11")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0857"
$c.Style = $savedStyle
$ws.Range("E11").Value = "  -4.16%  "
$c = $ws.Range("D12")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.865.97"
$c.Style = $savedStyle
$ws.Range("E12").Value = "  -0.25%  "
$c = $ws.Range("D13")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.640.29"
$c.Style = $savedStyle
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("E15").Value = "  -1.32%  "
$c = $ws.Range("D16")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "65.13"
$c.Style = $savedStyle
$c = $ws.Range("D17")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.667.30"
$c.Style = $savedStyle
$ws.Range("E17").Value = "  +0.04%  "
$c = $ws.Range("D18")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "230.25"
$c.Style = $savedStyle
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -0.32%  "
$c = $ws.Range("D20")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.60"
$c.Style = $savedStyle
$ws.Range("E20").Value = "  -1.43%  "
$c = $ws.Range("D21")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $savedStyle
$ws.Range("E21").Value = "  +0.00%  "
$c = $ws.Range("D22")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.61"
$c.Style = $savedStyle
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  +4.27%  "
$c = $ws.Range("D25")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "148.99"
$c.Style = $savedStyle
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -0.57%  "
$c = $ws.Range("D28")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "15.55"
$c.Style = $savedStyle
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.60%  "
$c = $ws.Range("D32")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = $savedStyle
$ws.Range("E32").Value = "  -0.95%  "
$c = $ws.Range("D33")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.483.45"
$c.Style = $savedStyle
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -1.76%  "
$c = $ws.Range("D36")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.33"
$c.Style = $savedStyle
$ws.Range("E36").Value = "  -1.34%  "
$c = $ws.Range("D37")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.963"
$c.Style = $savedStyle
$ws.Range("E37").Value = "  +7.79%  "
$c = $ws.Range("D38")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.882"
$c.Style = $savedStyle
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("E40").Value = "  +0.11%  "
$c = $ws.Range("D41")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.02"
$c.Style = $savedStyle
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -1.09%  "
$c = $ws.Range("D46")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.33"
$c.Style = $savedStyle
$ws.Range("E46").Value = "  -4.75%  "
$c = $ws.Range("D47")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.775.38"
$c.Style = $savedStyle
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  +0.44%  "
$c = $ws.Range("D49")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "87.57"
$c.Style = $savedStyle
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  -1.74%  "
$c = $ws.Range("D51")
$savedStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0990"
$c.Style = $savedStyle
$ws.Range("E51").Value = "  -0.30%  "
